# Apply "New format for EBE dream 11" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend the "score" defined name from A1:B7 to A1:C7
$wb.Names.Item("score").RefersTo = "=Sheet1!`$A`$1:`$C`$7"

# 2. Add the new "Fromat 2" column (column C) mirroring column B's layout/styles
$ws.Range("B1:B7").Copy() | Out-Null
$ws.Range("C1:C7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C1").Value = "Fromat 2"
$ws.Range("C2").Value = 40
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = -10
$ws.Range("C5").Value = -20
$ws.Range("C6").Value = -25
$ws.Range("C7").Value = 0

# 3. Update rows 44-52: set B column to 2, and update the VLOOKUP column index
#    in the D/G/J/M/P/S formulas from a hard-coded 2 to $B<row>+1
$map = @{
    "D" = "E"
    "G" = "H"
    "J" = "K"
    "M" = "N"
    "P" = "Q"
    "S" = "T"
}

for ($row = 44; $row -le 52; $row++) {
    $ws.Range("B$row").Value = 2

    foreach ($formulaCol in $map.Keys) {
        $lookupCol = $map[$formulaCol]
        $formula = '=IF(ISERROR(VLOOKUP(RANK(' + $lookupCol + $row + ', ($T' + $row + ',$Q' + $row + ',$N' + $row + ',$K' + $row + ',$H' + $row + ',$E' + $row + '), 0),  score, $B' + $row + '+1, FALSE)),"",VLOOKUP(RANK(' + $lookupCol + $row + ', ($T' + $row + ',$Q' + $row + ',$N' + $row + ',$K' + $row + ',$H' + $row + ',$E' + $row + '), 0),  score, $B' + $row + '+1, FALSE))'
        $ws.Range($formulaCol + $row).Formula = $formula
    }
}
